$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability matrix values (row sums remain 1)
# Row 2
$ws.Cells.Item(2, 2).Value = 0.2916666666666667
$ws.Cells.Item(2, 3).Value = 0.3958333333333333
$ws.Cells.Item(2, 10).Value = 0.02083333333333333
$ws.Cells.Item(2, 16).Value = 0.25
$ws.Cells.Item(2, 19).Value = 0.04166666666666666

# Row 3
$ws.Cells.Item(3, 10).Value = 0.1052631578947368
$ws.Cells.Item(3, 16).Value = 0.7368421052631579
$ws.Cells.Item(3, 19).Value = 0.1578947368421053

# Row 4
$ws.Cells.Item(4, 10).Value = 0.2
$ws.Cells.Item(4, 16).Value = 0.4
$ws.Cells.Item(4, 19).Value = 0.4

# Row 6
$ws.Cells.Item(6, 2).Value = 0.04761904761904762
$ws.Cells.Item(6, 10).Value = 0.3809523809523809
$ws.Cells.Item(6, 17).Value = 0.09523809523809523
$ws.Cells.Item(6, 18).Value = 0.04761904761904762
$ws.Cells.Item(6, 19).Value = 0.4285714285714285

# Row 7
$ws.Cells.Item(7, 2).Value = 0.25
$ws.Cells.Item(7, 10).Value = 0.25
$ws.Cells.Item(7, 17).Value = 0.375
$ws.Cells.Item(7, 19).Value = 0.125

# Row 8
$ws.Cells.Item(8, 2).Value = 0.0625
$ws.Cells.Item(8, 4).Value = 0.04166666666666666
$ws.Cells.Item(8, 6).Value = 0.04166666666666666
$ws.Cells.Item(8, 10).Value = 0.1458333333333333
$ws.Cells.Item(8, 15).Value = 0.0625
$ws.Cells.Item(8, 17).Value = 0.3333333333333333
$ws.Cells.Item(8, 18).Value = 0.0625
$ws.Cells.Item(8, 19).Value = 0.25

# Row 9
$ws.Cells.Item(9, 2).Value = 0.1071428571428571
$ws.Cells.Item(9, 4).Value = 0.03571428571428571
$ws.Cells.Item(9, 6).Value = 0.07142857142857142
$ws.Cells.Item(9, 10).Value = 0.1785714285714286
$ws.Cells.Item(9, 17).Value = 0.2142857142857143
$ws.Cells.Item(9, 18).Value = 0.07142857142857142
$ws.Cells.Item(9, 19).Value = 0.3214285714285715

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1159420289855072
$ws.Cells.Item(10, 4).Value = 0.00966183574879227
$ws.Cells.Item(10, 6).Value = 0.04347826086956522
$ws.Cells.Item(10, 10).Value = 0.1352657004830918
$ws.Cells.Item(10, 15).Value = 0.01449275362318841
$ws.Cells.Item(10, 17).Value = 0.3043478260869565
$ws.Cells.Item(10, 18).Value = 0.06280193236714976
$ws.Cells.Item(10, 19).Value = 0.3140096618357488

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1428571428571428
$ws.Cells.Item(11, 10).Value = 0.2857142857142857
$ws.Cells.Item(11, 11).Value = 0.1428571428571428
$ws.Cells.Item(11, 12).Value = 0.4285714285714285

# Row 12
$ws.Cells.Item(12, 7).Value = 0.6666666666666666
$ws.Cells.Item(12, 10).Value = 0.3333333333333333

# Row 13
$ws.Cells.Item(13, 7).Value = 0.8
$ws.Cells.Item(13, 10).Value = 0.2

# Row 14
$ws.Cells.Item(14, 7).Value = 1

# Row 15
$ws.Cells.Item(15, 8).Value = 0.1481481481481481
$ws.Cells.Item(15, 9).Value = 0.07407407407407407
$ws.Cells.Item(15, 10).Value = 0.4444444444444444
$ws.Cells.Item(15, 15).Value = 0.07407407407407407
$ws.Cells.Item(15, 19).Value = 0.2592592592592592

# Row 16
$ws.Cells.Item(16, 8).Value = 0.1481481481481481
$ws.Cells.Item(16, 10).Value = 0.3703703703703703
$ws.Cells.Item(16, 11).Value = 0.03703703703703703
$ws.Cells.Item(16, 13).Value = 0.03703703703703703
$ws.Cells.Item(16, 15).Value = 0.1851851851851852

# Row 17
$ws.Cells.Item(17, 6).Value = 0.01136363636363636
$ws.Cells.Item(17, 8).Value = 0.1704545454545454
$ws.Cells.Item(17, 9).Value = 0.1363636363636364
$ws.Cells.Item(17, 10).Value = 0.5454545454545454
$ws.Cells.Item(17, 11).Value = 0.02272727272727273
$ws.Cells.Item(17, 13).Value = 0.02272727272727273
$ws.Cells.Item(17, 14).Value = 0.01136363636363636
$ws.Cells.Item(17, 15).Value = 0.03409090909090909
$ws.Cells.Item(17, 19).Value = 0.04545454545454546

# Row 18
$ws.Cells.Item(18, 8).Value = 0.1111111111111111
$ws.Cells.Item(18, 9).Value = 0.05555555555555555
$ws.Cells.Item(18, 10).Value = 0.6666666666666666
$ws.Cells.Item(18, 11).Value = 0.05555555555555555
$ws.Cells.Item(18, 15).Value = 0.05555555555555555
$ws.Cells.Item(18, 19).Value = 0.05555555555555555

# Row 19
$ws.Cells.Item(19, 6).Value = 0.02362204724409449
$ws.Cells.Item(19, 8).Value = 0.1811023622047244
$ws.Cells.Item(19, 9).Value = 0.07874015748031496
$ws.Cells.Item(19, 10).Value = 0.5433070866141733
$ws.Cells.Item(19, 11).Value = 0.01574803149606299
$ws.Cells.Item(19, 13).Value = 0.01574803149606299
$ws.Cells.Item(19, 15).Value = 0.07086614173228346
$ws.Cells.Item(19, 19).Value = 0.07086614173228346
